$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember column A's width (the new column inherits A's formatting,
# matching the "width 17" shared by columns A and B in the target file).
$aWidth = $ws.Columns("A").ColumnWidth

# Insert a new column before the current column B ("Under 500gm"),
# shifting B:F right to C:G.
$ws.Columns("B").Insert()

# Give the freshly inserted column the same width as column A.
$ws.Columns("B").ColumnWidth = $aWidth

# The defined name "ListMonths" pointed at the old column C
# (Sheet1!$C$1, header "500-1000gm"); after the insert that header now
# lives in column D, so repoint the name to match.
$name = $wb.Names.Item("ListMonths")
$name.RefersTo = "=OFFSET(Sheet1!`$D`$1,1,0,MAX(Sheet1!`$A:`$A),1)"

# Populate the new "Free Shipping" column: header, plus a value of
# 1000 on the 282005 row (the only row that has one in the target).
$ws.Range("B1").Value = "Free Shipping"
$ws.Range("B4").Value = 1000

# Match the final selection left behind in the saved file.
$ws.Range("E8").Select() | Out-Null
